# Update automàtic: dades i banners [2026-02-06 23:19]
# Applies the per-cell text updates from the meteocat daily extraction refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding percentage-looking text (e.g. "74%") must keep their literal
# text type instead of being auto-parsed into a numeric percentage by Excel,
# so force the Text number format on them before writing the new value.
$percentCells = @('H3', 'H4', 'H8', 'H12', 'H22', 'H26', 'H31', 'H36')
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = '@'
}

$ws.Range('E2').Value = '2026-02-06 23:17:43'
$ws.Range('E3').Value = '2026-02-06 23:17:45'
$ws.Range('H3').Value = '74%'
$ws.Range('O3').Value = '-2.4 °C'
$ws.Range('E4').Value = '2026-02-06 23:17:47'
$ws.Range('H4').Value = '56%'
$ws.Range('J4').Value = '998.0 hPa'
$ws.Range('E5').Value = '2026-02-06 23:17:50'
$ws.Range('J5').Value = '998.2 hPa'
$ws.Range('O5').Value = '11.0 °C'
$ws.Range('E6').Value = '2026-02-06 23:17:52'
$ws.Range('J6').Value = '999.3 hPa'
$ws.Range('E7').Value = '2026-02-06 23:17:55'
$ws.Range('J7').Value = '999.0 hPa'
$ws.Range('E8').Value = '2026-02-06 23:17:57'
$ws.Range('H8').Value = '79%'
$ws.Range('O8').Value = '9.4 °C'
$ws.Range('E9').Value = '2026-02-06 23:18:00'
$ws.Range('O9').Value = '4.6 °C'
$ws.Range('E10').Value = '2026-02-06 23:18:02'
$ws.Range('O10').Value = '9.4 °C'
$ws.Range('E11').Value = '2026-02-06 23:18:05'
$ws.Range('J11').Value = '999.7 hPa'
$ws.Range('E12').Value = '2026-02-06 23:18:07'
$ws.Range('H12').Value = '65%'
$ws.Range('O12').Value = '13.0 °C'
$ws.Range('E13').Value = '2026-02-06 23:18:09'
$ws.Range('E14').Value = '2026-02-06 23:18:12'
$ws.Range('N14').Value = '-6.3 °C 22:55 TU'
$ws.Range('E15').Value = '2026-02-06 23:18:14'
$ws.Range('J15').Value = '998.4 hPa'
$ws.Range('K15').Value = '11.6 MJ/m2'
$ws.Range('O15').Value = '10.1 °C'
$ws.Range('E16').Value = '2026-02-06 23:18:17'
$ws.Range('E17').Value = '2026-02-06 23:18:19'
$ws.Range('I17').Value = '0.3 mm'
$ws.Range('J17').Value = '999.7 hPa'
$ws.Range('E18').Value = '2026-02-06 23:18:22'
$ws.Range('E19').Value = '2026-02-06 23:18:24'
$ws.Range('I19').Value = '3.5 mm'
$ws.Range('J19').Value = '1000.6 hPa'
$ws.Range('O19').Value = '9.6 °C'
$ws.Range('E20').Value = '2026-02-06 23:18:27'
$ws.Range('E21').Value = '2026-02-06 23:18:29'
$ws.Range('J21').Value = '998.6 hPa'
$ws.Range('E22').Value = '2026-02-06 23:18:32'
$ws.Range('H22').Value = '81%'
$ws.Range('O22').Value = '9.9 °C'
$ws.Range('E23').Value = '2026-02-06 23:18:34'
$ws.Range('J23').Value = '998.4 hPa'
$ws.Range('O23').Value = '9.8 °C'
$ws.Range('E24').Value = '2026-02-06 23:18:36'
$ws.Range('J24').Value = '997.8 hPa'
$ws.Range('E25').Value = '2026-02-06 23:18:39'
$ws.Range('J25').Value = '999.3 hPa'
$ws.Range('E26').Value = '2026-02-06 23:18:41'
$ws.Range('H26').Value = '79%'
$ws.Range('E27').Value = '2026-02-06 23:18:44'
$ws.Range('O27').Value = '11.0 °C'
$ws.Range('E28').Value = '2026-02-06 23:18:46'
$ws.Range('J28').Value = '1000.6 hPa'
$ws.Range('E29').Value = '2026-02-06 23:18:49'
$ws.Range('O29').Value = '12.0 °C'
$ws.Range('E30').Value = '2026-02-06 23:18:51'
$ws.Range('E31').Value = '2026-02-06 23:18:54'
$ws.Range('H31').Value = '85%'
$ws.Range('I31').Value = '4.1 mm'
$ws.Range('J31').Value = '1000.1 hPa'
$ws.Range('E32').Value = '2026-02-06 23:18:56'
$ws.Range('K32').Value = '12.0 MJ/m2'
$ws.Range('O32').Value = '15.1 °C'
$ws.Range('E33').Value = '2026-02-06 23:18:58'
$ws.Range('O33').Value = '10.0 °C'
$ws.Range('E34').Value = '2026-02-06 23:19:00'
$ws.Range('E35').Value = '2026-02-06 23:19:03'
$ws.Range('N35').Value = '-4.1 °C 22:57 TU'
$ws.Range('E36').Value = '2026-02-06 23:19:06'
$ws.Range('H36').Value = '70%'
$ws.Range('J36').Value = '1000.7 hPa'
$ws.Range('N36').Value = '6.0 °C 22:52 TU'
$ws.Range('O36').Value = '12.2 °C'
